# Atualizei dados da bibi e add
# - Insere um novo registro diario (Dia 4, Junho/2025) na planilha de faturamento.
# - Atualiza os totais de venda (B) de alguns dias ja existentes (Junho e Maio/2025).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere uma nova linha na posicao 4 (empurra os registros de Maio/Abril/Marco uma linha para baixo)
# e preenche os dados do novo dia: Dia=4, total_venda=5167.9, Mes=6, Ano=2025, Periodo="06/2025".
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 5167.9
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 2025
$ws.Range("E4").Value = "06/2025"

# Atualiza os totais de venda de Junho/2025 ja existentes.
$ws.Range("B2").Value = 42933.45
$ws.Range("B3").Value = 18600.98

# Atualiza os totais de venda de Maio/2025 ja existentes (linhas deslocadas em +1 apos o insert acima).
$ws.Range("B17").Value = 21827.07
$ws.Range("B18").Value = 9475.47
$ws.Range("B23").Value = 24291.06
$ws.Range("B24").Value = 30467.22
$ws.Range("B25").Value = 21933.56
